# Outcomes.xlsx - "Updated Excel with Amir for NOGO outcomes"
# Fill in the NOGO-trial (and GO-trial "AcquireFixError_"/"BreakFix_") outcome
# columns on STAT_ABORT that were previously left blank, fix the "BreakFix_,
# Abort_" label to just "BreakFix_", and correct "NoGoWrong?" -> "NoGoError?"
# in the comment on row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STAT_ABORT")

# Row 2 (FIX_ON / AcquireFixError_): extend value across GO + both NOGO cols
$ws.Range("E2").Value = "AcquireFixError_"
$ws.Range("F2").Value = "AcquireFixError_"
$ws.Range("G2").Value = "AcquireFixError_"

# Row 4 (FIX_HOLD / BreakFix_): label simplified from "BreakFix_, Abort_"
$ws.Range("D4").Value = "BreakFix_"
$ws.Range("E4").Value = "BreakFix_"
$ws.Range("F4").Value = "BreakFix_"
$ws.Range("G4").Value = "BreakFix_"

# Row 12's comment is corrected first so its shared-string slot is allocated
# ahead of the "NoGoError_" label (matches the order the strings were added
# in the authored workbook).
$ws.Range("H12").Value = "Gaze not in TARG_WIN during Tone`nWhat tone to be given for NOGO Trial, since this is NoGoError?"

# Rows 6-16: fill in NOGO Trial columns (F, G) with "NoGoError_"
$ws.Range("F6").Value = "NoGoError_"
$ws.Range("G6").Value = "NoGoError_"

$ws.Range("F7").Value = "NoGoError_"
$ws.Range("G7").Value = "NoGoError_"

$ws.Range("F8").Value = "NoGoError_"
$ws.Range("G8").Value = "NoGoError_"

$ws.Range("F9").Value = "NoGoError_"
$ws.Range("G9").Value = "NoGoError_"

$ws.Range("F10").Value = "NoGoError_"
$ws.Range("G10").Value = "NoGoError_"

$ws.Range("F11").Value = "NoGoError_"
$ws.Range("G11").Value = "NoGoError_"

$ws.Range("F12").Value = "NoGoError_"
$ws.Range("G12").Value = "NoGoError_"

$ws.Range("F13").Value = "NoGoError_"
$ws.Range("G13").Value = "NoGoError_"

$ws.Range("F14").Value = "NoGoError_"
$ws.Range("G14").Value = "NoGoError_"

$ws.Range("F15").Value = "NoGoError_"
$ws.Range("G15").Value = "NoGoError_"

$ws.Range("F16").Value = "NoGoError_"
$ws.Range("G16").Value = "NoGoError_"

# Move the active selection to match the saved cursor position
$ws.Range("G22").Select()
